# Adds a new "2022-Q4" quarterly sheet right after "总计", updates the
# summary ("总计") sheet with the new quarter's aggregate row, and shifts
# the existing quarterly data down by one position (names/content unchanged).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: Update the "总计" (summary) sheet.
#   Shift existing data rows 2-6 down to rows 3-7 (keeps all formatting),
#   then re-write column A as a simple 0-based row index, and fill in the
#   new row 2 with the 2022-Q4 totals.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

for ($r = 6; $r -ge 2; $r--) {
  $src = $ws1.Range("A" + $r + ":D" + $r)
  $dstRow = $r + 1
  $dst = $ws1.Range("A" + $dstRow + ":D" + $dstRow)
  $src.Copy($dst)
}

for ($r = 2; $r -le 7; $r++) {
  $ws1.Range("A" + $r).Value2 = $r - 2
}

$ws1.Range("B2").Value2 = "2022-Q4"
$ws1.Range("C2").Value2 = 3
$ws1.Range("D2").Value2 = 0.1

# ---------------------------------------------------------------------------
# Step 2: Insert a brand-new worksheet named "2022-Q4" immediately after
# "总计" (i.e. before the sheet currently in position 2, "2022-Q3"), so the
# final sheet order becomes: 总计, 2022-Q4, 2022-Q3, 2022-Q1, 2021-Q4,
# 2021-Q3, 2020-Q4.
# ---------------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($wb.Worksheets.Item(2))
$newSheet.Name = "2022-Q4"

# Match the sheetPr / page-setup conventions used by the rest of the
# workbook (outline summary rows below, summary columns to the right;
# the same 0.75/0.75/1/1/0.5/0.5 inch margins as every other sheet).
$newSheet.Outline.SummaryRow = 1
$newSheet.Outline.SummaryColumn = 1
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

# Borrow the header-row / index-column formatting (bold font, border,
# centered alignment) from one of the existing quarterly sheets so the new
# sheet matches the workbook's established look - that sheet is now at
# position 3 ("2022-Q3") after the insert above.
$refSheet = $wb.Worksheets.Item(3)
$refSheet.Range("B1:H1").Copy($newSheet.Range("B1:H1"))
$refSheet.Range("A2").Copy($newSheet.Range("A2"))
$refSheet.Range("A2").Copy($newSheet.Range("A3"))
$refSheet.Range("A2").Copy($newSheet.Range("A4"))

# Header labels.
$newSheet.Range("B1").Value2 = "基金代码"
$newSheet.Range("C1").Value2 = "基金名称"
$newSheet.Range("D1").Value2 = "基金规模"
$newSheet.Range("E1").Value2 = "股票总仓位"
$newSheet.Range("F1").Value2 = "仓位占比"
$newSheet.Range("G1").Value2 = "持有市值(亿元)"
$newSheet.Range("H1").Value2 = "仓位排名"

# Row 2 - 光大保德信量化股票
$newSheet.Range("A2").Value2 = 0
$newSheet.Range("B2").Value2 = "'360001"
$newSheet.Range("C2").Value2 = "光大保德信量化股票"
$newSheet.Range("D2").Value2 = "'12.20"
$newSheet.Range("E2").Value2 = "'88.99"
$newSheet.Range("F2").Value2 = "'0.76"
$newSheet.Range("G2").Value2 = "'0.0927"
$newSheet.Range("H2").Value2 = 3

# Row 3 - 光大保德信锦弘混合A
$newSheet.Range("A3").Value2 = 1
$newSheet.Range("B3").Value2 = "'011231"
$newSheet.Range("C3").Value2 = "光大保德信锦弘混合A"
$newSheet.Range("D3").Value2 = "'1.25"
$newSheet.Range("E3").Value2 = "'29.26"
$newSheet.Range("F3").Value2 = "'0.31"
$newSheet.Range("G3").Value2 = "'0.0039"
$newSheet.Range("H3").Value2 = 4

# Row 4 - 光大保德信锦弘混合C
$newSheet.Range("A4").Value2 = 2
$newSheet.Range("B4").Value2 = "'011232"
$newSheet.Range("C4").Value2 = "光大保德信锦弘混合C"
$newSheet.Range("D4").Value2 = "'0.71"
$newSheet.Range("E4").Value2 = "'29.26"
$newSheet.Range("F4").Value2 = "'0.31"
$newSheet.Range("G4").Value2 = "'0.0022"
$newSheet.Range("H4").Value2 = 4

# The leading "'" above forces Excel to keep numeric-looking text (fund
# codes / decimal figures) as literal strings instead of auto-converting
# them to numbers, but it also stamps the cells with a Text number format.
# Strip that back off (without touching the values) by pasting in the
# (blank, unformatted) number format from an untouched cell.
$newSheet.Range("Z100").Copy()
$newSheet.Range("B2:G4").PasteSpecial(-4122)

# Restore the originally-active tab (the last sheet, "2020-Q4") since
# adding a new sheet makes it the active/selected one by default.
$wb.Worksheets.Item($wb.Worksheets.Count).Activate()

Write-Host "2022-Q4 sheet added and 总计 updated"
